$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Helper: append one new data row (date / running-total formula / daily delta)
# to the bottom of one of the three data sheets, copying the formatting of
# the previous last row (row 84) down onto the new row (row 85).
# ---------------------------------------------------------------------------
function Add-DataRow {
    param($SheetName, $LastRow, $DateSerial, $DailyValue)

    $ws = $wb.Worksheets.Item($SheetName)
    $newRow = $LastRow + 1

    # Copy the formatting (style) of the last existing row down to the new row
    $ws.Range("A" + $LastRow + ":C" + $LastRow).Copy() | Out-Null
    $ws.Range("A" + $newRow + ":C" + $newRow).PasteSpecial(-4122) | Out-Null

    # A column: date
    $ws.Cells.Item($newRow, 1).Value = $DateSerial

    # B column: running total = previous running total + this row's C value
    $ws.Cells.Item($newRow, 2).Formula = "=SUM(B" + $LastRow + "+C" + $newRow + ")"

    # C column: the day's raw value
    $ws.Cells.Item($newRow, 3).Value = $DailyValue
}

Add-DataRow "Confirmed" 84 43981 1764
Add-DataRow "Recoverd"  84 43981 360
Add-DataRow "Death"     84 43981 28

# ---------------------------------------------------------------------------
# Update the selection shown on each sheet to track the newly-added last row,
# the same way Excel leaves the selection where the user last worked.
# ---------------------------------------------------------------------------
$wsConfirmed = $wb.Worksheets.Item("Confirmed")
$wsConfirmed.Range("B84:B85").Select() | Out-Null

$wsRecoverd = $wb.Worksheets.Item("Recoverd")
$wsRecoverd.Range("B84:B85").Select() | Out-Null

$wsDeath = $wb.Worksheets.Item("Death")
$wsDeath.Range("B84:B85").Select() | Out-Null

# ---------------------------------------------------------------------------
# The active tab moves from "Death" (index 2) to "Recoverd" (index 1).
# Activating a sheet makes Excel mark it as the selected / active tab and
# updates the workbook's stored active-tab index accordingly.
# ---------------------------------------------------------------------------
$wsRecoverd.Activate()
$wsRecoverd.Range("B84:B85").Select() | Out-Null
